$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$mdUrl57 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8616aea046e367c2a199740543ec507c8bec236/e2e/57441d4c-9c7c-4ef2-b274-607766f11a1b.md"
$mdUrl82 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8616aea046e367c2a199740543ec507c8bec236/e2e/82504e87-6a4a-430b-abba-4adbe191f342.md"

# --- Status column (C) text update: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime column (K) update ---
$wsZhCn.Range("K2").Value = "2016-08-27 16:29:46"
$wsZhCn.Range("K3").Value = "2016-08-27 16:29:46"
$wsDeDe.Range("K2").Value = "2016-08-27 16:29:53"
$wsDeDe.Range("K3").Value = "2016-08-27 16:29:53"

# --- Latest Target File (I) + Latest Handback File (J) for zh-cn ---
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl57, "", "", "57441d4c-9c7c-4ef2-b274-607766f11a1b.md")
$wsZhCn.Range("J2").Value = "57441d4c-9c7c-4ef2-b274-607766f11a1b.a8391b9233de920d614aa5f3225aac6ab30777c6.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl82, "", "", "82504e87-6a4a-430b-abba-4adbe191f342.md")
$wsZhCn.Range("J3").Value = "82504e87-6a4a-430b-abba-4adbe191f342.d3504e8e32041521409db48875af57a9f636eb2f.zh-cn.xlf"

# --- Latest Target File (I) + Latest Handback File (J) for de-de ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl57, "", "", "57441d4c-9c7c-4ef2-b274-607766f11a1b.md")
$wsDeDe.Range("J2").Value = "57441d4c-9c7c-4ef2-b274-607766f11a1b.a8391b9233de920d614aa5f3225aac6ab30777c6.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl82, "", "", "82504e87-6a4a-430b-abba-4adbe191f342.md")
$wsDeDe.Range("J3").Value = "82504e87-6a4a-430b-abba-4adbe191f342.d3504e8e32041521409db48875af57a9f636eb2f.de-de.xlf"

# --- Column width adjustments ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "done"
